# chore: adapt column header formatting to respective input file names
#
# Renames the "_old" / "_new" column-header suffixes to the respective
# input-file-version suffixes ("_FV2410" / "_FV2504"), turns the header
# row + data range into a native Excel Table ("Table1"), and freezes the
# header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells -------------------------------------------------
# Columns A-J currently end in "_old"  -> should end in "_FV2410"
# Columns L-U currently end in "_new"  -> should end in "_FV2504"
# Column K ("diff") is left untouched.

$newHeadersFV2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

$newHeadersFV2504 = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

# Columns 1..10 => A..J
for ($i = 0; $i -lt $newHeadersFV2410.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeadersFV2410[$i]
}

# Columns 12..21 => L..U (column 11 / K = "diff" is unchanged)
for ($i = 0; $i -lt $newHeadersFV2504.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeadersFV2504[$i]
}

# --- 2. Turn the range into an Excel Table ("Table1") ---------------------------
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))

$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row -----------------------------------------------------
$ws.Cells.Item(2, 1).Activate()
$excel.ActiveWindow.FreezePanes = $true
